$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10: replace "RatioConstraints" (A10 / B10=0.1) with "RatiosBlastFurnace" (A10 / C10=0.2)
$ws.Range("A10").Value = "RatiosBlastFurnace"
$ws.Range("B10").ClearContents()
$ws.Range("C10").Value = 0.2

# Update the view: active cell A10, no fixed scroll position
$ws.Range("A10").Select()
